# Update ratings employee model
# The user changed the Attitude-vs-Work pairwise comparison on the
# "pairwise_comp" sheet: cell F7 now holds the direct judgement (4,
# "Attitude strongly preferred over Work") instead of its reciprocal
# (0.25); its mirror cell E8 now holds the reciprocal (0.25) instead of
# the direct value (4). This re-derives all of the dependent
# (pre-computed / hard-coded) priority values across the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "pairwise_comp" : Criteria pairwise-comparison matrix
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("pairwise_comp")

# F7 becomes the "input" judgement cell (style changes from the
# locked/computed look, border2, to the editable upper-triangle look,
# border1) and its value changes 0.25 -> 4.
$ws1.Range("F7").Value = 4
foreach ($idx in 7, 8, 9, 10) {
    $ws1.Range("F7").Borders.Item($idx).LineStyle = 1
    $ws1.Range("F7").Borders.Item($idx).Weight = 2
    $ws1.Range("F7").Borders.Item($idx).ColorIndex = -4105
}

# E8 keeps its existing (lower-triangle / reciprocal) style, only the
# value changes 4 -> 0.25.
$ws1.Range("E8").Value = 0.25

# Recomputed "Normal" (K) / "Ideal" (L) priority columns, and the
# "Incons." cell, rows 4-10.
$ws1.Range("K4").Value = 0.077
$ws1.Range("L4").Value = 0.213
$ws1.Range("K5").Value = 0.193
$ws1.Range("L5").Value = 0.531
$ws1.Range("K6").Value = 0.049
$ws1.Range("L6").Value = 0.135
$ws1.Range("K7").Value = 0.363
$ws1.Range("K8").Value = 0.083
$ws1.Range("L8").Value = 0.23
$ws1.Range("L9").Value = 0.649
$ws1.Range("L10").Value = 0.102

# ---------------------------------------------------------------
# Sheet "supermatrix" : Criteria column (C4:C8) of the supermatrix
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("supermatrix")
$ws3.Range("C4").Value = 0.077
$ws3.Range("C5").Value = 0.193
$ws3.Range("C6").Value = 0.049
$ws3.Range("C7").Value = 0.363
$ws3.Range("C8").Value = 0.083

# ---------------------------------------------------------------
# Sheet "limit matrix" : Criteria column (C4:C11) of the limit matrix
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("limit matrix")
$ws4.Range("C4").Value = 0.07099999999999999
$ws4.Range("C5").Value = 0.178
$ws4.Range("C6").Value = 0.045
$ws4.Range("C7").Value = 0.335
$ws4.Range("C8").Value = 0.077
$ws4.Range("C9").Value = 0.217
$ws4.Range("C10").Value = 0.047
$ws4.Range("C11").Value = 0.03

# ---------------------------------------------------------------
# Sheet "limitingPriorities" : priorities column (C3:C10)
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("limitingPriorities")
$ws5.Range("C3").Value = 0.07099999999999999
$ws5.Range("C4").Value = 0.178
$ws5.Range("C5").Value = 0.045
$ws5.Range("C6").Value = 0.335
$ws5.Range("C7").Value = 0.077
$ws5.Range("C8").Value = 0.217
$ws5.Range("C9").Value = 0.047
$ws5.Range("C10").Value = 0.03

# ---------------------------------------------------------------
# Sheet "localPriorities" : priorities column (C3:C7)
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("localPriorities")
$ws6.Range("C3").Value = 0.077
$ws6.Range("C4").Value = 0.193
$ws6.Range("C5").Value = 0.049
$ws6.Range("C6").Value = 0.363
$ws6.Range("C7").Value = 0.083
